$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Prix Spot")

# New header for column BA (05-aug)
$ws.Range("BA1").Value = "05-aug"

# New numeric values for BA2:BA25
$values = @(
    8.460000000000001,
    -0.07000000000000001,
    -0.01,
    -0.06,
    -0.25,
    -0.11,
    0.02,
    2.77,
    0,
    -0.03,
    -2.78,
    -3.54,
    -9.65,
    -24.02,
    -19.5,
    -14.94,
    -14.44,
    -0.02,
    3,
    46.4,
    73.5,
    82.56,
    88.22,
    74.95
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Range("BA$row").Value = $values[$i]
}

# Match header style used by the rest of row 1 (bold/centered/bordered)
$ws.Range("AZ1").Copy()
$ws.Range("BA1").PasteSpecial(-4122)
